{"js": "const pairs = [\n  [\"2025-05-10 Saturday\", \"2025-05-11 Sunday\"],\n  [\"47\u00d788=4136\", \"94\u00d724=2256\"],\n  [\"72\u00d740=2880\", \"32\u00d753=1696\"],\n  [\"86\u00d754=4644\", \"92\u00d731=2852\"],\n  [\"99\u00d789=8811\", \"25\u00d750=1250\"],\n  [\"68\u00d763=4284\", \"28\u00d744=1232\"],\n  [\"66\u00d759=3894\", \"64\u00d725=1600\"],\n  [\"18\u00d740=720\", \"93\u00d750=4650\"],\n  [\"77\u00d732=2464\", \"19\u00d728=532\"],\n  [\"81\u00d719=1539\", \"31\u00d756=1736\"],\n  [\"59\u00d757=3363\", \"68\u00d737=2516\"],\n  [\"24\u00d716=384\", \"39\u00d719=741\"],\n  [\"49\u00d718=882\", \"23\u00d714=322\"],\n  [\"86\u00d724=2064\", \"32\u00d765=2080\"],\n  [\"87\u00d711=957\", \"86\u00d771=6106\"],\n  [\"89\u00d796=8544\", \"33\u00d776=2508\"],\n  [\"21\u00d791=1911\", \"15\u00d727=405\"],\n  [\"67\u00d741=2747\", \"52\u00d762=3224\"],\n  [\"45\u00d788=3960\", \"47\u00d728=1316\"],\n  [\"46\u00d732=1472\", \"93\u00d757=5301\"],\n  [\"35\u00d774=2590\", \"85\u00d783=7055\"],\n  [\"99\u00d788=8712\", \"55\u00d757=3135\"],\n  [\"39\u00d762=2418\", \"35\u00d736=1260\"],\n  [\"14\u00d762=868\", \"34\u00d746=1564\"],\n  [\"83\u00d758=4814\", \"44\u00d769=3036\"],\n  [\"95\u00d765=6175\", \"46\u00d713=598\"]\n];\n\nconst body = context.document.body;\nlet totalFound = 0;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  totalFound += results.items.length;\n}\nawait context.sync();\nreturn totalFound;\n", "ps1": "$pairs = @(\n    @('2025-05-10 Saturday', '2025-05-11 Sunday'),\n    @('47\u00d788=4136', '94\u00d724=2256'),\n    @('72\u00d740=2880', '32\u00d753=1696'),\n    @('86\u00d754=4644', '92\u00d731=2852'),\n    @('99\u00d789=8811', '25\u00d750=1250'),\n    @('68\u00d763=4284', '28\u00d744=1232'),\n    @('66\u00d759=3894', '64\u00d725=1600'),\n    @('18\u00d740=720', '93\u00d750=4650'),\n    @('77\u00d732=2464', '19\u00d728=532'),\n    @('81\u00d719=1539', '31\u00d756=1736'),\n    @('59\u00d757=3363', '68\u00d737=2516'),\n    @('24\u00d716=384', '39\u00d719=741'),\n    @('49\u00d718=882', '23\u00d714=322'),\n    @('86\u00d724=2064', '32\u00d765=2080'),\n    @('87\u00d711=957', '86\u00d771=6106'),\n    @('89\u00d796=8544', '33\u00d776=2508'),\n    @('21\u00d791=1911', '15\u00d727=405'),\n    @('67\u00d741=2747', '52\u00d762=3224'),\n    @('45\u00d788=3960', '47\u00d728=1316'),\n    @('46\u00d732=1472', '93\u00d757=5301'),\n    @('35\u00d774=2590', '85\u00d783=7055'),\n    @('99\u00d788=8712', '55\u00d757=3135'),\n    @('39\u00d762=2418', '35\u00d736=1260'),\n    @('14\u00d762=868', '34\u00d746=1564'),\n    @('83\u00d758=4814', '44\u00d769=3036'),\n    @('95\u00d765=6175', '46\u00d713=598')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found for replacement: $oldText\"\n    }\n}\n\nWrite-Output \"Replaced $($pairs.Count) items\"\n"}
